# Insert a new data row at row 287 (pushes existing rows 287:370 down to 288:371)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(287).Insert()

# Populate the newly inserted row 287 with the new record's values
$ws.Range("A287").Value = 5
$ws.Range("B287").Value = "Macroferia Regional de Talca"
$ws.Range("C287").Value = "Maule"
$ws.Range("D287").Value = 44841
$ws.Range("E287").Value = 7
$ws.Range("F287").Value = 100112003
$ws.Range("G287").Value = "Ajo"
$ws.Range("H287").Value = "Chino"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 300
$ws.Range("K287").Value = 20000
$ws.Range("L287").Value = 20000
$ws.Range("M287").Value = 20000
$ws.Range("N287").Value = "`$/malla 10 kilos"
$ws.Range("O287").Value = "China"
$ws.Range("P287").Value = 2000
$ws.Range("Q287").Value = 10
$ws.Range("R287").Value = "Hortaliza"
